# Roll back the "batch uploading examinees" feature changes:
# the sheet previously gained a 3-row block (rows 1-3) holding
# "科目列表1/2/3" sample subject lists plus an extra duplicated
# examinee table. Revert to the original 4-row layout and restore
# the original subject-list text in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra "科目列表" rows that were inserted above the table.
$ws.Rows("1:3").Delete()

# Restore the original subject-list values for the remaining rows.
$ws.Range("E2").Value = "语文,数学,英语"
$ws.Range("E3").Value = "语文,数学,英语,物理,化学"
$ws.Range("E4").Value = "语文,数学(文),英语,政治,历史"

# Restore the original selected cell.
$ws.Range("E20").Select()

# Remove the duplicated Hyperlink / Followed Hyperlink cell styles that
# had accumulated (keep only the first of each, plus Normal).
$styles = $wb.Styles
$toDelete = @(10, 9, 8, 7, 5, 4, 3, 2)
foreach ($i in $toDelete) {
    $styles.Item($i).Delete()
}
